$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update T3:W3 values from 0.99 to 1
$ws.Range("T3:W3").Value = 1

# Update the selection range shown in the sheet view
$ws.Range("A1:X7").Select()
